$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing quarterly
# columns (old D:K) right to F:M.
$ws.Range("D:E").EntireColumn.Insert()

# The newly inserted D:E columns come in with default formatting; copy
# number formats from column F (which now holds the data that used to be
# in column D) so D:E match the rest of the quarterly columns.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarterly columns (D = newest quarter, E = next).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 537800
$ws.Range("E8").Value = 514400
$ws.Range("D9").Value = 462900
$ws.Range("E9").Value = 448600
$ws.Range("D10").Value = 74900
$ws.Range("E10").Value = 65800
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 12200
$ws.Range("E15").Value = 11900
$ws.Range("D17").Value = 503300
$ws.Range("E17").Value = 485100
$ws.Range("D18").Value = 34500
$ws.Range("E18").Value = 29300
$ws.Range("D20").Value = 600
$ws.Range("E20").Value = 400
$ws.Range("D21").Value = 47200
$ws.Range("E21").Value = 41700
$ws.Range("D22").Value = 3700
$ws.Range("E22").Value = 4000
$ws.Range("D23").Value = 31300
$ws.Range("E23").Value = 25800
$ws.Range("D24").Value = 4800
$ws.Range("E24").Value = 5400
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 26600
$ws.Range("E26").Value = 20400
$ws.Range("D27").Value = 26400
$ws.Range("E27").Value = 20900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -600
$ws.Range("E32").Value = -400
$ws.Range("D33").Value = 26400
$ws.Range("E33").Value = 20900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 26400
$ws.Range("E35").Value = 20900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 31100
$ws.Range("E41").Value = 45700
$ws.Range("D42").Value = 8700
$ws.Range("E42").Value = 9700
$ws.Range("D43").Value = 276100
$ws.Range("E43").Value = 261500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 32200
$ws.Range("E45").Value = 31600
$ws.Range("D46").Value = 348100
$ws.Range("E46").Value = 348400
$ws.Range("D47").Value = 36300
$ws.Range("E47").Value = 32600
$ws.Range("D48").Value = 618900
$ws.Range("E48").Value = 593100
$ws.Range("D49").Value = 139100
$ws.Range("E49").Value = 136400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 39700
$ws.Range("E52").Value = 33100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1182000
$ws.Range("E54").Value = 1143600
$ws.Range("D57").Value = 44200
$ws.Range("E57").Value = 41300
$ws.Range("D58").Value = 10100
$ws.Range("E58").Value = 10100
$ws.Range("D59").Value = 214900
$ws.Range("E59").Value = 196500
$ws.Range("D60").Value = 269200
$ws.Range("E60").Value = 247900
$ws.Range("D61").Value = 233100
$ws.Range("E61").Value = 245600
$ws.Range("D62").Value = 77300
$ws.Range("E62").Value = 77700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 591000
$ws.Range("E66").Value = 581700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 344900
$ws.Range("E72").Value = 321400
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 590900
$ws.Range("E76").Value = 561900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 26400
$ws.Range("E81").Value = 20900
$ws.Range("D83").Value = 12200
$ws.Range("E83").Value = 11900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 53000
$ws.Range("E89").Value = 56000
$ws.Range("D91").Value = -17300
$ws.Range("E91").Value = -13300
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -55900
$ws.Range("E94").Value = -14000
$ws.Range("D96").Value = -2400
$ws.Range("E96").Value = -2300
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -11700
$ws.Range("E100").Value = -23500
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -14600
$ws.Range("E102").Value = 18500

# Row 91 ("Capital Expenditures") values for the five oldest-shown
# quarters were restated; fix F:J after the shift (K:M already match the
# values that shifted over correctly).
$ws.Range("F91").Value = -19800
$ws.Range("G91").Value = -4400
$ws.Range("H91").Value = 100
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
